# Fill in the homework grades that were still recorded as 0 with the
# actual scores collected for this grading pass.
#
# For every student row (4-30), any of the ДЗ_1..ДЗ_4 cells (columns C:F)
# that is currently 0 gets bumped to 2, except for Назаралиев Расул's
# ДЗ_1 (C20) which is updated to 5 and loses the "not submitted" green
# highlight (style changes to the plain bordered style used by other
# completed cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 4; $row -le 30; $row++) {
    foreach ($col in @("C", "D", "E", "F")) {
        $cell = $ws.Range("$col$row")
        if ($cell.Value2 -eq 0) {
            if ($col -eq "C" -and $row -eq 20) {
                # This particular mark turned out to be a full score (5), so
                # besides updating the number it also loses the "missing
                # homework" green highlight - match the plain bordered look
                # already used by the other graded cells (e.g. C21).
                $ws.Range("C21").Copy()
                $cell.PasteSpecial(-4122)  # xlPasteFormats
                $cell.Value2 = 5
            } else {
                $cell.Value2 = 2
            }
        }
    }
}

$excel.CutCopyMode = 0

$ws.Range("C4").Select()
